# Update "想去人数" (number of people interested) values per the latest
# scrape output (commit: "Update gh-pages to output generated at 456a3b4").
#
# Changes are split across three worksheets:
#   - 展览 (Exhibitions)
#   - 演出 (Performances)
#   - 全部类型 (All types, an aggregate of the sheets above)

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4678
$ws1.Range("F3").Value  = 2721
$ws1.Range("F5").Value  = 2736
$ws1.Range("F9").Value  = 1714
$ws1.Range("F12").Value = 208
$ws1.Range("F13").Value = 395
$ws1.Range("F16").Value = 89
$ws1.Range("F22").Value = 642
$ws1.Range("F27").Value = 1654
$ws1.Range("F28").Value = 1455
$ws1.Range("F29").Value = 314
$ws1.Range("F31").Value = 1420
$ws1.Range("F32").Value = 2278
$ws1.Range("F33").Value = 380
$ws1.Range("F36").Value = 113
$ws1.Range("F39").Value = 769
$ws1.Range("F40").Value = 1455
$ws1.Range("F43").Value = 483
$ws1.Range("F44").Value = 19
$ws1.Range("F46").Value = 106

# --- 演出 sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 87
$ws2.Range("F11").Value = 31

# --- 全部类型 sheet ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4678
$ws4.Range("F3").Value  = 2721
$ws4.Range("F4").Value  = 2736
$ws4.Range("F5").Value  = 1714
$ws4.Range("F10").Value = 208
$ws4.Range("F11").Value = 395
$ws4.Range("F14").Value = 89
$ws4.Range("F19").Value = 642
$ws4.Range("F22").Value = 87
$ws4.Range("F27").Value = 1654
$ws4.Range("F28").Value = 1455
$ws4.Range("F29").Value = 314
$ws4.Range("F33").Value = 2278
$ws4.Range("F34").Value = 380
$ws4.Range("F38").Value = 31
$ws4.Range("F40").Value = 113
$ws4.Range("F43").Value = 769
$ws4.Range("F44").Value = 1455
$ws4.Range("F47").Value = 483
$ws4.Range("F49").Value = 106

$wb.Save()
